$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextValue 'D2' '96.425.97'
$ws.Range('E2').Value = '  -0.96%  '

Set-TextValue 'D3' '3.637.46'
$ws.Range('E3').Value = '  -2.75%  '

Set-TextValue 'D4' '2.72'
$ws.Range('E4').Value = '  +43.20%  '

Set-TextValue 'D5' '1.01'
$ws.Range('E5').Value = '  +0.69%  '

Set-TextValue 'D6' '225.81'
$ws.Range('E6').Value = '  -5.53%  '

Set-TextValue 'D7' '640.34'
$ws.Range('E7').Value = '  -3.29%  '

Set-TextValue 'D8' '0.425'
$ws.Range('E8').Value = '  -4.03%  '

Set-TextValue 'D9' '1.20'
$ws.Range('E9').Value = '  +12.08%  '

Set-TextValue 'D10' '1.00'
$ws.Range('E10').Value = '  +0.30%  '

Set-TextValue 'D11' '3.633.94'
$ws.Range('E11').Value = '  -2.81%  '

Set-TextValue 'D12' '48.32'
$ws.Range('E12').Value = '  +6.91%  '

$ws.Range('E13').Value = '  +1.50%  '

Set-TextValue 'D14' '0.0000291'
$ws.Range('E14').Value = '  -9.57%  '

Set-TextValue 'D15' '6.50'
$ws.Range('E15').Value = '  -6.58%  '

Set-TextValue 'D16' '4.321.92'
$ws.Range('E16').Value = '  -2.58%  '

Set-TextValue 'D17' '96.222.33'
$ws.Range('E17').Value = '  -0.93%  '

Set-TextValue 'D18' '21.30'
$ws.Range('E18').Value = '  +12.70%  '

Set-TextValue 'D19' '8.86'
$ws.Range('E19').Value = '  -2.56%  '

Set-TextValue 'D20' '14.02'
$ws.Range('E20').Value = '  +6.37%  '

Set-TextValue 'D21' '3.632.92'
$ws.Range('E21').Value = '  -2.66%  '

Set-TextValue 'D22' '0.566'
$ws.Range('E22').Value = '  +11.44%  '

Set-TextValue 'D23' '0.280'
$ws.Range('E23').Value = '  +46.26%  '

Set-TextValue 'D24' '515.99'
$ws.Range('E24').Value = '  -2.81%  '

Set-TextValue 'D25' '3.25'
$ws.Range('E25').Value = '  -7.64%  '

Set-TextValue 'D26' '120.39'
$ws.Range('E26').Value = '  +11.47%  '

Set-TextValue 'D27' '0.0000200'
$ws.Range('E27').Value = '  -12.82%  '

Set-TextValue 'D28' '6.78'
$ws.Range('E28').Value = '  -2.02%  '

Set-TextValue 'D29' '3.828.75'
$ws.Range('E29').Value = '  -2.63%  '

Set-TextValue 'D30' '12.77'
$ws.Range('E30').Value = '  -6.75%  '

Set-TextValue 'D31' '12.84'
$ws.Range('E31').Value = '  -0.25%  '

Set-TextValue 'D32' '3.01'
$ws.Range('E32').Value = '  -1.57%  '

Set-TextValue 'D33' '1.00'
$ws.Range('E33').Value = '  +0.14%  '

$ws.Range('B34').Value = 'PolygonEcosystemToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue 'D34' '0.619'
$ws.Range('E34').Value = '  +3.73%  '

Set-TextValue 'D35' '32.73'
$ws.Range('E35').Value = '  +0.07%  '

$ws.Range('B36').Value = 'Cronos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D36' '0.180'
$ws.Range('E36').Value = '  -6.45%  '

$ws.Range('E37').Value = '  +0.63%  '

Set-TextValue 'D38' '1.75'
$ws.Range('E38').Value = '  -5.35%  '

$ws.Range('B39').Value = 'USDe'
$ws.Range('C39').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D39' '1.00'
$ws.Range('E39').Value = '  +0.00%  '

$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue 'D40' '8.34'
$ws.Range('E40').Value = '  -5.41%  '

$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D41' '582.56'
$ws.Range('E41').Value = '  -10.30%  '

Set-TextValue 'D42' '6.99'
$ws.Range('E42').Value = '  +2.79%  '

$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D43' '0.496'
$ws.Range('E43').Value = '  +3.50%  '

$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D44' '41.17'
$ws.Range('E44').Value = '  +0.89%  '

$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D45' '0.0502'
$ws.Range('E45').Value = '  +9.02%  '

Set-TextValue 'D46' '0.157'
$ws.Range('E46').Value = '  -5.42%  '

Set-TextValue 'D47' '0.954'
$ws.Range('E47').Value = '  -2.77%  '

Set-TextValue 'D48' '1.93'
$ws.Range('E48').Value = '  -5.92%  '

Set-TextValue 'D49' '231.10'
$ws.Range('E49').Value = '  +11.45%  '

$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D50' '8.80'
$ws.Range('E50').Value = '  +0.76%  '

$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D51' '2.27'
$ws.Range('E51').Value = '  -5.16%  '
